$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: reorder set contents 'any','list','List[any]' -> 'list','any','List[any]'
$ws.Range("E2").Value = "{'list', 'any', 'List[any]'}"

# E3: 'any' -> 'list'
$ws.Range("E3").Value = "list"

# F3: was a "Loss" (red fill) -> now "Neutral" (orange fill, same style as F2/F4)
$ws.Range("F3").Value = "Neutral"
$ws.Range("F3").Interior.Color = $ws.Range("F2").Interior.Color()

# D5: PyType Wins count 1 -> 0
$ws.Range("D5").Value = 0

# Insert a new row above the old row 6 (pushes old row 6 down to row 7,
# carrying its formatting/content with it)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the Scalpel accuracy summary
$ws.Range("C6").Value = "Scalpel Accuracy:"
$ws.Range("D6").Value = 100

# Update the (now shifted) accuracy-over-pytype value on row 7
$ws.Range("F7").Value = 100
